# CDS_SPARSE_2020_2021.xlsx edit
# - Redefine "High School Units" entities: drop units-required/units-recommended,
#   add units/require/recommend (one "units" column always = 1, plus separate
#   require/recommend indicator columns split out of the old M/N columns).
# - Make "High School Units" the active/selected sheet (was Enrollment_General).

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("High School Units")

# --- Header row: rename/introduce the three entity columns ---
$ws3.Range("M1").Value = "units"
$ws3.Range("N1").Value = "require"
$ws3.Range("O1").Value = "recommend"

# --- Data rows 2-22: units is always required (1); the old M/N values move
#     to the new N (require) / O (recommend) columns respectively ---
for ($r = 2; $r -le 22; $r++) {
    $oldM = $ws3.Cells.Item($r, 13).Value()
    $oldN = $ws3.Cells.Item($r, 14).Value()
    $ws3.Cells.Item($r, 13).Value = 1
    $ws3.Cells.Item($r, 14).Value = $oldM
    $ws3.Cells.Item($r, 15).Value = $oldN
}

# --- Column width for the newly introduced column O ---
$ws3.Range("O1").ColumnWidth = 13.2

# --- Make "High School Units" the active sheet/tab (was Enrollment_General) ---
$ws3.Activate()
$ws3.Range("N25").Select()
